$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.44174971086602
$ws.Range("C2").Value = 0.2640045443680492
$ws.Range("E2").Value = 0.09596278397741465
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.002431580731880127
$ws.Range("L2").Value = 0.219633493926338
$ws.Range("O2").Value = 2.390920195325748
$ws.Range("B3").Value = 1.303140114108146
$ws.Range("C3").Value = 0.2491817021483769
$ws.Range("E3").Value = 0.09721320375492404
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.002434612479328385
$ws.Range("L3").Value = 0.2093033950525154
$ws.Range("O3").Value = 2.431778622624023
$ws.Range("B4").Value = 1.218024777860876
$ws.Range("C4").Value = 0.2400666556826536
$ws.Range("E4").Value = 0.09804550853818128
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.002436571690183954
$ws.Range("L4").Value = 0.2030581434057979
$ws.Range("O4").Value = 2.459274812000729
$ws.Range("B5").Value = 1.183339318764922
$ws.Range("C5").Value = 0.2363490481393455
$ws.Range("E5").Value = 0.09840090219977249
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002437394728222686
$ws.Range("L5").Value = 0.2005377278387783
$ws.Range("O5").Value = 2.471083789701794
$ws.Range("B6").Value = 1.177579863928827
$ws.Range("C6").Value = 0.23573156101574
$ws.Range("E6").Value = 0.09846089465513153
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002437532883847697
$ws.Range("L6").Value = 0.2001207009989514
$ws.Range("O6").Value = 2.473081095364051
$ws.Range("B7").Value = 1.217556996225653
$ws.Range("C7").Value = 0.240016531077373
$ws.Range("E7").Value = 0.09805023582543448
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.00243658269036956
$ws.Range("L7").Value = 0.2030240526072475
$ws.Range("O7").Value = 2.459431628454055
$ws.Range("B8").Value = 1.39395995063694
$ws.Range("C8").Value = 0.2588966636075156
$ws.Range("E8").Value = 0.0963805289373969
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.002432605844709387
$ws.Range("L8").Value = 0.2160514809068417
$ws.Range("O8").Value = 2.40450718872053
$ws.Range("B9").Value = 1.73975455471782
$ws.Range("C9").Value = 0.2957997796023335
$ws.Range("E9").Value = 0.09361865697893457
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002425579088038114
$ws.Range("L9").Value = 0.2423703286894181
$ws.Range("O9").Value = 2.315990132728004
$ws.Range("B10").Value = 1.993671093130786
$ws.Range("C10").Value = 0.3228262143030349
$ws.Range("E10").Value = 0.0919023120619773
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002420882240507374
$ws.Range("L10").Value = 0.2621780942683074
$ws.Range("O10").Value = 2.262758908652387
$ws.Range("B11").Value = 2.1091437457697
$ws.Range("C11").Value = 0.3351000648873708
$ws.Range("E11").Value = 0.09118951568940048
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002418845620880235
$ws.Range("L11").Value = 0.2712917998551205
$ws.Range("O11").Value = 2.241127654197072
$ws.Range("B12").Value = 2.152863750121583
$ws.Range("C12").Value = 0.3397446319501967
$ws.Range("E12").Value = 0.09092938184653754
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.00241808870821739
$ws.Range("L12").Value = 0.2747577177325553
$ws.Range("O12").Value = 2.233309867180935
$ws.Range("B13").Value = 2.143448205687093
$ws.Range("C13").Value = 0.3387444918209042
$ws.Range("E13").Value = 0.09098497087505741
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002418251087593061
$ws.Range("L13").Value = 0.2740106150552037
$ws.Range("O13").Value = 2.234976925158435
$ws.Range("B14").Value = 2.112740767656931
$ws.Range("C14").Value = 0.3354822435848916
$ws.Range("E14").Value = 0.09116791811578651
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.00241878306289057
$ws.Range("L14").Value = 0.2715766476447783
$ws.Range("O14").Value = 2.240476983572734
$ws.Range("B15").Value = 2.093930609467577
$ws.Range("C15").Value = 0.3334835868823518
$ws.Range("E15").Value = 0.09128125344040861
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002419110774552677
$ws.Range("L15").Value = 0.2700876920446404
$ws.Range("O15").Value = 2.243894626317342
$ws.Range("B16").Value = 1.986123810104459
$ws.Range("C16").Value = 0.3220236479608332
$ws.Range("E16").Value = 0.09195026385448202
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002421017344749403
$ws.Range("L16").Value = 0.2615845604096592
$ws.Range("O16").Value = 2.264224722368766
$ws.Range("B17").Value = 1.919977430852327
$ws.Range("C17").Value = 0.314987845751233
$ws.Range("E17").Value = 0.09237809900899663
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.002422212527016537
$ws.Range("L17").Value = 0.2563945134032082
$ws.Range("O17").Value = 2.277359755777837
$ws.Range("B18").Value = 1.881928577363169
$ws.Range("C18").Value = 0.3109391167244553
$ws.Range("E18").Value = 0.09263057615183357
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002422909380745618
$ws.Range("L18").Value = 0.2534190392630222
$ws.Range("O18").Value = 2.28515777241121
$ws.Range("B19").Value = 1.869045397523848
$ws.Range("C19").Value = 0.3095679661151678
$ws.Range("E19").Value = 0.09271715902792543
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.002423146942975808
$ws.Range("L19").Value = 0.2524132634103751
$ws.Range("O19").Value = 2.287839739187604
$ws.Range("B20").Value = 1.927019171567338
$ws.Range("C20").Value = 0.3157370200570995
$ws.Range("E20").Value = 0.09233189300666922
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002422084323819896
$ws.Range("L20").Value = 0.2569459989445875
$ws.Range("O20").Value = 2.275936335792608
$ws.Range("B21").Value = 2.121760491200689
$ws.Range("C21").Value = 0.336440536266565
$ws.Range("E21").Value = 0.09111391637682864
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.0024186264208174
$ws.Range("L21").Value = 0.2722911629177389
$ws.Range("O21").Value = 2.238851330985852
$ws.Range("B22").Value = 2.248993629461609
$ws.Range("C22").Value = 0.3499522583138344
$ws.Range("E22").Value = 0.09037494781642152
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002416449866247155
$ws.Range("L22").Value = 0.2824060915957745
$ws.Range("O22").Value = 2.216792075719979
$ws.Range("B23").Value = 2.181091338076669
$ws.Range("C23").Value = 0.3427426580650774
$ws.Range("E23").Value = 0.09076412538550471
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002417603927283056
$ws.Range("L23").Value = 0.2769997175871026
$ws.Range("O23").Value = 2.228365594037371
$ws.Range("B24").Value = 1.923835664051467
$ws.Range("C24").Value = 0.315398330020713
$ws.Range("E24").Value = 0.09235276246623592
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002422142254096518
$ws.Range("L24").Value = 0.2566966463105302
$ws.Range("O24").Value = 2.276579096235167
$ws.Range("B25").Value = 1.646228174407383
$ws.Range("C25").Value = 0.2858307651922019
$ws.Range("E25").Value = 0.09431091522747792
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002427397882168111
$ws.Range("L25").Value = 0.2351676659169186
$ws.Range("O25").Value = 2.337871227421431
